$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 224 (existing rows 224:293 shift down to 225:294),
# which also extends the used range / dimension to R294.
$ws.Rows.Item(224).Insert()

# Populate the newly inserted row 224 with the new record (Ajo / Chino, 60 units,
# 2022-11-11 date serial 44876), matching the price/unit fields already used for
# the "Chilote" variety batch that used to occupy this row position.
$ws.Cells.Item(224, 1).Value  = 7
$ws.Cells.Item(224, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(224, 3).Value  = "Ñuble"
$ws.Cells.Item(224, 4).Value  = 44876
$ws.Cells.Item(224, 5).Value  = 16
$ws.Cells.Item(224, 6).Value  = 100112003
$ws.Cells.Item(224, 7).Value  = "Ajo"
$ws.Cells.Item(224, 8).Value  = "Chino"
$ws.Cells.Item(224, 9).Value  = "Primera"
$ws.Cells.Item(224, 10).Value = 60
$ws.Cells.Item(224, 11).Value = 14000
$ws.Cells.Item(224, 12).Value = 15000
$ws.Cells.Item(224, 13).Value = 14500
$ws.Cells.Item(224, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(224, 15).Value = "China"
$ws.Cells.Item(224, 16).Value = 1450
$ws.Cells.Item(224, 17).Value = 10
$ws.Cells.Item(224, 18).Value = "Hortaliza"
